$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The DATE column (B) for every logged row (rows 2-80) needs to move from
# 2024-09-16 to 2024-09-17, now that the model training/validation run has
# completed. The values are plain text (not real Excel dates), so we build
# the replacement through a literal-text formula and then convert it back
# to a static value via copy / paste-special. This avoids Excel's automatic
# "looks like a date" reinterpretation (which would turn the text into a
# date serial number and change the cell's number format/style).
for ($row = 2; $row -le 80; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Formula = "=""2024-09-17"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}
